# "added support for groups of size 5"
#
# 1) Add a new "Params" sheet (holds the shared weighting parameter) in front
#    of the existing M2/M3/M4 sheets.
# 2) Point every module sheet's "Prop" cell at Params!B1 via a formula
#    instead of a hard-coded literal.
# 3) Fix M4 (groups of 4), which was missing its "total points" row, by
#    adding the SUM row back in (same pattern as M2/M3).
# 4) Add a new "M5" sheet (groups of 5) at the end, following the same
#    layout as M2/M3/M4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Params sheet, inserted before the first existing sheet (M2).
# ---------------------------------------------------------------------
$wsParams = $wb.Worksheets.Add()
$wsParams.Name = "Params"
$wsParams.Range("A1").Value = "Prop"
$wsParams.Range("B1").Value = 0.5
$wsParams.Range("B2").Select()

# ---------------------------------------------------------------------
# 2) Wire up M2 / M3 / M4 "Prop" cells to reference Params!B1.
# ---------------------------------------------------------------------
$wsM2 = $wb.Worksheets.Item("M2")
$wsM2.Range("F1").Formula = "=Params!B1"
$wsM2.Range("B4").Select()

$wsM3 = $wb.Worksheets.Item("M3")
$wsM3.Range("G1").Formula = "=Params!B1"
$wsM3.Range("B5").Select()

$wsM4 = $wb.Worksheets.Item("M4")
$wsM4.Range("H1").Formula = "=Params!B1"

# ---------------------------------------------------------------------
# 3) M4 was missing the totals row (row 6) under the four module columns
#    -- add it back, matching the SUM pattern used on the other sheets.
# ---------------------------------------------------------------------
$wsM4.Range("B6").Formula = "=SUM(B2:B5)"
$wsM4.Range("C6:E6").Formula = "=SUM(C2:C5)"
$wsM4.Range("H1").Select()

# ---------------------------------------------------------------------
# 4) New M5 sheet (groups of 5), appended after M4.
# ---------------------------------------------------------------------
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsM5 = $wb.Worksheets.Add($null, $wsLast)
$wsM5.Name = "M5"

$wsM5.Range("B1").Value = "M1"
$wsM5.Range("C1").Value = "M2"
$wsM5.Range("D1").Value = "M3"
$wsM5.Range("E1").Value = "M4"
$wsM5.Range("F1").Value = "M5"
$wsM5.Range("H1").Value = "Prop"
$wsM5.Range("I1").Formula = "=Params!B1"

$wsM5.Range("A2").Value = "M1"
$wsM5.Range("H2").Value = "Group Grade"
$wsM5.Range("I2").Value = 100

$wsM5.Range("A3").Value = "M2"
$wsM5.Range("A4").Value = "M3"
$wsM5.Range("A5").Value = "M4"
$wsM5.Range("A6").Value = "M5"

$wsM5.Range("B7").Formula = "=SUM(B2:B6)"
$wsM5.Range("C7:F7").Formula = "=SUM(C2:C6)"

$wsM5.Range("B8").Formula = "=B7/SUM(`$B7:`$F7)"
$wsM5.Range("C8:F8").Formula = "=C7/SUM(`$B7:`$F7)"

$wsM5.Range("A11").Value = "Individual Grades"
$wsM5.Range("B11").Formula = "=(`$I`$1*`$I`$2*B8*5)+(`$I`$2*(1-`$I`$1))"
$wsM5.Range("C11:F11").Formula = "=(`$I`$1*`$I`$2*C8*5)+(`$I`$2*(1-`$I`$1))"

$wsM5.Range("F11").Select()
